$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 3 (FAPs/Cort/Ghsr/FAPs ...) entirely - its data is being
# folded into row 2, and the second data row disappears from the sheet.
$ws.Rows("3:3").Delete()

# Row 2 now carries what used to be row 3's values (Target cluster switches
# from "ECs" to "FAPs"), except the four "derived specificity" columns
# (O, P, S, T) which become 1 instead of the non-normalized figure.
$ws.Range("D2").Value = "FAPs"
$ws.Range("M2").Value = 0.1546876666666667
$ws.Range("N2").Value = 0.464063
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.01885106406088889
$ws.Range("R2").Value = 0.169659576548
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
